$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style from B13 (already-filled row) to B14:B15
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14:B15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 14: 16-Jan-2021, 1PM-4PM, Core Java, Download and Installation of JDK...
$ws.Range("B14").Value = 44212
$ws.Range("C14").Value = "1PM-4PM"
$ws.Range("D14").Value = "Core Java"
$ws.Range("E14").Value = "Download and Installation of JDK, added path and coded basic programs"

# Row 15: 16-Jan-2021, 7PM-9PM, Core Java, Started with OOPS concept in JAVA
$ws.Range("B15").Value = 44212
$ws.Range("C15").Value = "7PM-9PM"
$ws.Range("D15").Value = "Core Java"
$ws.Range("E15").Value = "Started with OOPS concept in JAVA"

$ws.Range("E15").Select() | Out-Null
